$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1) values to the new lowercase/pluralized names,
# keeping the underlying data (row 2) unchanged.
$ws.Range("A1").Value = "hgnc_id"
$ws.Range("B1").Value = "hgnc_symbol"
$ws.Range("C1").Value = "disease_associated_transcripts"
$ws.Range("D1").Value = "genetic_disease_models"
$ws.Range("E1").Value = "mosaicism"
$ws.Range("F1").Value = "reduced_penetrance"
$ws.Range("G1").Value = "database_entry_version"

# Update the selected cell to match the new active selection (G1).
$ws.Range("G1").Select()
